$d = $word.ActiveDocument
Write-Host "Paragraphs count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host $i ":" $p.Range.Text
}
